$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

$c = $ws.Cells.Item(162, 2)
$c.Value = "07 12:19>>> 0FD8A6BA00   John Orellana"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(163, 2)
$c.Value = "07 12:21>>> CFD89A9C80   Bryan Williams"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(164, 2)
$c.Value = "07 12:21>>> 0FD8AF9E00   Karla Pimentel"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(165, 2)
$c.Value = "07 12:21>>> 0FD8AF9E00   Karla Pimentel"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(166, 2)
$c.Value = "07 12:43>>> 0FD8B4E900   Carlos De Los Santos"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(167, 2)
$c.Value = "07 13:07>>> 0FD8AFD720   Jennifer Mercedes"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(168, 2)
$c.Value = "07 13:10>>> 4FD8A85BA0   hojin euam"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(169, 2)
$c.Value = "07 13:11>>> 8FD8985540    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(170, 2)
$c.Value = "07 13:32>>> 8FD8A8D440    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(171, 2)
$c.Value = "07 13:34>>> 8FD8ADBD20   david orlando"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(172, 2)
$c.Value = "07 13:42>>> 8FD8B39BA0   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(173, 2)
$c.Value = "07 13:53>>> 8FD8A91340   Jesse Silkworth"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(174, 2)
$c.Value = "07 14:07>>> CFD8A948E0   Vanessa Cordero"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(175, 2)
$c.Value = "07 15:28>>> 8FD8AFD100   rahyner penaranda"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(176, 2)
$c.Value = "07 15:40>>> 0FD8A6BA00   John Orellana"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(177, 2)
$c.Value = "07 15:48>>> 0FD8AD42A0   Covalky Pena"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(178, 2)
$c.Value = "07 16:12>>> 0FD8A9BD80   Jaspreet Kaur"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(179, 2)
$c.Value = "07 16:23>>> CFD8A3DFA0   eugene marmontov"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(180, 2)
$c.Value = "07 17:19>>> 0FD88F2580   Carolyn Alana"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(181, 2)
$c.Value = "07 17:47>>> 4FD87F2960   Thomas Yoo"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(182, 2)
$c.Value = "07 17:48>>> CFD893A460   Gary Tsai"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(183, 2)
$c.Value = "07 19:10>>> 0FD8AD42A0   Covalky Pena"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(184, 2)
$c.Value = "07 19:10>>> 0FD8B5ED00   Crosby anne"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(185, 2)
$c.Value = "07 19:55>>> CFD893A460   Gary Tsai"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(186, 2)
$c.Value = "08 08:53>>> 4FD8A33DE0   Anne Crosby"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(187, 2)
$c.Value = "08 08:54>>> 0FD8AE8B60   Natalie Primus"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(188, 2)
$c.Value = "08 09:07>>> 8FD8AEB240   luciano.ibbott"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(189, 2)
$c.Value = "08 09:14>>> 8FD8B68DE0   Miguel Martillo"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(190, 2)
$c.Value = "08 09:42>>> 4FD8A51080   Calvin Y Au"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(191, 2)
$c.Value = "08 09:52>>> 8FD8AC1E00   Kimberly Pierre"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(192, 2)
$c.Value = "08 09:54>>> 4FD8B41A40   Justin Davis"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(193, 2)
$c.Value = "08 10:03>>> 0FD89CD000   David Schachner"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(194, 2)
$c.Value = "08 10:04>>> CFD8B45C40   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(195, 2)
$c.Value = "08 10:04>>> CFD8A6B7E0   peter sormilic"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(196, 2)
$c.Value = "08 10:28>>> 0FD8B5ED00   Crosby anne"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(197, 2)
$c.Value = "08 10:30>>> 0FD8A83600   douglas smith"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(198, 2)
$c.Value = "08 10:37>>> 8FD88BA9C0    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(199, 2)
$c.Value = "08 11:11>>> 0FD899D9C0   cammy-el allen"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(200, 2)
$c.Value = "08 11:26>>> 0FD8ACF700   Jared Amuso"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(201, 2)
$c.Value = "08 11:30>>> 0FD89E2980   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(202, 2)
$c.Value = "08 11:41>>> 8FD8AEA680   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(203, 2)
$c.Value = "08 11:58>>> 4FD8A33DE0   Anne Crosby"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(204, 2)
$c.Value = "08 12:21>>> 4FD87D1F40   Gao Feng"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(205, 2)
$c.Value = "08 12:22>>> 0FD88F2580   Carolyn Alana"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(206, 2)
$c.Value = "08 12:49>>> CFD8B3CA60   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(207, 2)
$c.Value = "08 12:49>>> 4FD8AB6C20   Freddy Beltran"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(208, 2)
$c.Value = "08 13:09>>> 0FD8AE84A0    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(209, 2)
$c.Value = "08 13:18>>> CFD8B6A840   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(210, 2)
$c.Value = "08 13:23>>> 4FD8B3E300   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(211, 2)
$c.Value = "08 13:25>>> 8FD8AFF180   Oscar Chicaiza"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(212, 2)
$c.Value = "08 13:25>>> 4FD8B41A40   Justin Davis"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(213, 2)
$c.Value = "08 13:27>>> 8FD8A5D640   Anwar Sufian"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(214, 2)
$c.Value = "08 13:33>>> 8FD8A3A820   cowell,truman"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(215, 2)
$c.Value = "08 13:35>>> 8FD87CA8C0   Dwayne Stallworth"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(216, 2)
$c.Value = "08 13:36>>> 8FD8B367A0    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(217, 2)
$c.Value = "08 13:36>>> 0FD8AE0B20   Anthony Rivera"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(218, 2)
$c.Value = "08 13:38>>> 0FD8AC1480   "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(219, 2)
$c.Value = "08 13:39>>> 0FD8AF13C0   Rondell Holland"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(220, 2)
$c.Value = "08 13:40>>> 8FD8A57200   Averill Curameng"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(221, 2)
$c.Value = "08 15:02>>> 0FD8B5ED00   Crosby anne"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(222, 2)
$c.Value = "08 15:44>>> 0FD8AD42A0   Covalky Pena"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(223, 2)
$c.Value = "08 15:45>>> 0FD8A87380   Crystal Evelyn "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(224, 2)
$c.Value = "08 15:50>>> CFD8AFA4C0   arango juan"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(225, 2)
$c.Value = "08 15:52>>> 4FD8A85BA0   hojin euam"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(226, 2)
$c.Value = "08 16:11>>> 8FD8A5E5E0   Anthony Avevor"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(227, 2)
$c.Value = "08 16:27>>> 0FD89CD000   David Schachner"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(228, 2)
$c.Value = "08 17:14>>> 8FD8A5D640   Anwar Sufian"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(229, 2)
$c.Value = "08 18:04>>> 0FD8B5ED00   Crosby anne"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(230, 2)
$c.Value = "08 18:52>>> CFD8A78940   Nicole Latta"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(231, 2)
$c.Value = "08 19:11>>> CFD8B45C40    "
$c.NumberFormat = "General"
$c = $ws.Cells.Item(232, 2)
$c.Value = "08 20:41>>> CFD893A460   Gary Tsai"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(233, 2)
$c.Value = "08 20:44>>> 4FD8A85BA0   hojin euam"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(234, 2)
$c.Value = "08 20:46>>> 4FD8A85BA0   hojin euam"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(235, 2)
$c.Value = "08 20:47>>> 0FD8AF13C0   Rondell Holland"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(236, 2)
$c.Value = "08 20:51>>> 8FD894EAC0   Sade Thomas"
$c.NumberFormat = "General"
$c = $ws.Cells.Item(237, 2)
$c.Value = "08 20:51>>> 8FD894EAC0   Sade Thomas"
$c.NumberFormat = "General"
